# Regenerate the "K" column (G) values in save_data sheet.
# (Commit: regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2..38 (one per game row), replacing the old Strike# derived values.
$kValues = @(
    0, 1, 2, 1, 1, 2, 1, 0, 0, 0, 1, 2, 2, 1, 0, 1, 1, 1, 1, 2,
    1, 0, 0, 2, 1, 0, 1, 0, 0, 0, 3, 2, 1, 1, 0, 1, 0
)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
